# Energy Consumption4.xlsx update
# - Updates B2:C49 with new computed values
# - Removes row 50 (shrinks data from 49 rows to 48 rows), which also
#   changes the sheet dimension from A1:C50 to A1:C49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; B=5.991447168658981; C=2.491189000677316},
    @{Row=3; B=7.553647287603898; C=5.308026905412556},
    @{Row=4; B=8.97351535407549; C=7.970697773261305},
    @{Row=5; B=9.502765294463623; C=10.59610199492633},
    @{Row=6; B=9.76890494358604; C=13.05406209436645},
    @{Row=7; B=11.03875545935195; C=15.75176644168219},
    @{Row=8; B=16.51028764101113; C=18.38571699399124},
    @{Row=9; B=18.48159174730587; C=21.05375390917795},
    @{Row=10; B=20.37008767784808; C=23.59667721424405},
    @{Row=11; B=20.94986247447342; C=26.46998769500571},
    @{Row=12; B=23.61217926183835; C=29.39756608037219},
    @{Row=13; B=23.790546128254; C=31.84550543824598},
    @{Row=14; B=24.36752059986437; C=34.91670648996631},
    @{Row=15; B=25.46758679449792; C=37.74307846494585},
    @{Row=16; B=28.62412791623088; C=40.37325119781437},
    @{Row=17; B=30.57779685158517; C=43.12755947402447},
    @{Row=18; B=34.111869239285; C=46.11002434557073},
    @{Row=19; B=35.06856897145691; C=48.58350262917232},
    @{Row=20; B=37.53826485290286; C=51.15165500868702},
    @{Row=21; B=40.02571042475367; C=53.8294136002698},
    @{Row=22; B=42.57091060949232; C=56.58114307130248},
    @{Row=23; B=45.75081368089247; C=59.22775605894903},
    @{Row=24; B=48.68522483378318; C=62.4107562106937},
    @{Row=25; B=49.10015834519454; C=65.13662112326833},
    @{Row=26; B=52.05214139030311; C=67.69692320125336},
    @{Row=27; B=53.33118676499563; C=70.27287967232822},
    @{Row=28; B=54.55045537017131; C=72.54610971052402},
    @{Row=29; B=55.27318176760544; C=75.4535483700629},
    @{Row=30; B=57.21155732935492; C=78.07614666271212},
    @{Row=31; B=58.78544254228267; C=80.73297759496231},
    @{Row=32; B=59.08653849260632; C=83.48865774221173},
    @{Row=33; B=60.55700120896246; C=86.60187388146996},
    @{Row=34; B=61.08942347169119; C=89.32237183076558},
    @{Row=35; B=63.70488118465479; C=92.00638021722325},
    @{Row=36; B=65.93817406507419; C=94.70032574078766},
    @{Row=37; B=67.83808843508643; C=97.33147533118648},
    @{Row=38; B=70.15170793608013; C=100.0859080319831},
    @{Row=39; B=72.43047020906285; C=102.6741104127252},
    @{Row=40; B=74.65459814332687; C=105.329827278619},
    @{Row=41; B=75.71330306555116; C=107.9570459288677},
    @{Row=42; B=78.54441231709049; C=110.7301012896943},
    @{Row=43; B=79.6293050329584; C=113.6781633172621},
    @{Row=44; B=80.92008531908057; C=116.228814297003},
    @{Row=45; B=81.84334921639241; C=119.0716907024139},
    @{Row=46; B=92.13755480645804; C=121.6480316794211},
    @{Row=47; B=92.75385600099858; C=124.2307982821843},
    @{Row=48; B=93.60446988005721; C=126.9469451914802},
    @{Row=49; B=96.61530090904623; C=129.6359283558295}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
}

# Remove the now-obsolete last row (old row 50, A=48) so the table ends at row 49.
$ws.Rows.Item(50).Delete()
